$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H40").Value = 7229.846
$ws.Range("I40").Value = 4999.5
$ws.Range("K40").Value = 4999.5
$ws.Range("M40").Value = -4824.5
$ws.Range("H53").Value = 456.51852
$ws.Range("I53").Value = 806.8
$ws.Range("K53").Value = 806.8
$ws.Range("M53").Value = -169.8
$ws.Range("H98").Value = 2488
$ws.Range("I98").Value = 2550.2812
$ws.Range("J98").Value = 495
$ws.Range("K98").Value = 2550.2812
$ws.Range("L98").Value = 495
$ws.Range("M98").Value = -1052.2812
$ws.Range("N98").Value = -3491
$ws.Range("H100").Value = 4987.25
$ws.Range("I100").Value = 4987.25
$ws.Range("K100").Value = 4987.25
$ws.Range("M100").Value = -4446.25
$ws.Range("H101").Value = 3616.8
$ws.Range("I101").Value = 441.55554
$ws.Range("J101").Value = 8379.666999999999
$ws.Range("K101").Value = 1324.66662
$ws.Range("L101").Value = 25139.001
$ws.Range("M101").Value = 297.33338
$ws.Range("N101").Value = -28383.001
$ws.Range("H122").Value = 2488
$ws.Range("I122").Value = 2550.2812
$ws.Range("J122").Value = 495
$ws.Range("K122").Value = 7650.8436
$ws.Range("L122").Value = 1485
$ws.Range("M122").Value = -5200.8436
$ws.Range("N122").Value = -6385
$ws.Range("H132").Value = 1663.5385
$ws.Range("I132").Value = 1377.5
$ws.Range("K132").Value = 4132.5
$ws.Range("M132").Value = -1602.5
$ws.Range("H138").Value = 2705.1
$ws.Range("J138").Value = 3797.6924
$ws.Range("L138").Value = 11393.0772
$ws.Range("N138").Value = -21673.0772
$ws.Range("H141").Value = 3447.1785
$ws.Range("I141").Value = 3404.8462
$ws.Range("K141").Value = 10214.5386
$ws.Range("M141").Value = -5034.5386

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 5411.8076
$ws.Range("I2").Value = 4913.8667
$ws.Range("J2").Value = 6090.8184
$ws.Range("K2").Value = 4913.8667
$ws.Range("L2").Value = 6090.8184
$ws.Range("M2").Value = -4800.8667
$ws.Range("N2").Value = -6316.8184
$ws.Range("H8").Value = 2833.3333
$ws.Range("I8").Value = 1750
$ws.Range("K8").Value = 1750
$ws.Range("M8").Value = -1606
$ws.Range("H32").Value = 4664.061
$ws.Range("I32").Value = 3798.2444
$ws.Range("K32").Value = 3798.2444
$ws.Range("M32").Value = -3511.2444
$ws.Range("H110").Value = 3908.35
$ws.Range("I110").Value = 2246.1853
$ws.Range("K110").Value = 2246.1853
$ws.Range("M110").Value = -201.1853000000001
$ws.Range("H116").Value = 5411.8076
$ws.Range("I116").Value = 4913.8667
$ws.Range("J116").Value = 6090.8184
$ws.Range("K116").Value = 4913.8667
$ws.Range("L116").Value = 6090.8184
$ws.Range("M116").Value = -2619.8667
$ws.Range("N116").Value = -10678.8184

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 5411.8076
$ws.Range("I3").Value = 4913.8667
$ws.Range("J3").Value = 6090.8184
$ws.Range("K3").Value = 4913.8667
$ws.Range("L3").Value = 6090.8184
$ws.Range("M3").Value = -4799.8667
$ws.Range("N3").Value = -6318.8184
$ws.Range("H99").Value = 4488.3335
$ws.Range("I99").Value = 2709.2307
$ws.Range("K99").Value = 2709.2307
$ws.Range("M99").Value = -1211.2307

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 4149.048
$ws.Range("I31").Value = 3349.9546
$ws.Range("K31").Value = 3349.9546
$ws.Range("M31").Value = -3054.9546
$ws.Range("H34").Value = 4149.048
$ws.Range("I34").Value = 3349.9546
$ws.Range("K34").Value = 3349.9546
$ws.Range("M34").Value = -3147.9546
$ws.Range("H41").Value = 5309.5713
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H60").Value = 10216.667
$ws.Range("I60").Value = 6500
$ws.Range("K60").Value = 6500
$ws.Range("M60").Value = -5989
$ws.Range("H86").Value = 14318.667
$ws.Range("I86").Value = 13978.167
$ws.Range("J86").Value = 14999.667
$ws.Range("K86").Value = 13978.167
$ws.Range("L86").Value = 14999.667
$ws.Range("M86").Value = -12855.167
$ws.Range("N86").Value = -17245.667
$ws.Range("H89").Value = 14318.667
$ws.Range("I89").Value = 13978.167
$ws.Range("J89").Value = 14999.667
$ws.Range("K89").Value = 69890.83499999999
$ws.Range("L89").Value = 74998.33499999999
$ws.Range("M89").Value = -64274.83499999999
$ws.Range("N89").Value = -86230.33499999999
$ws.Range("H97").Value = 28419.666
$ws.Range("J97").Value = 28103.6
$ws.Range("L97").Value = 28103.6
$ws.Range("N97").Value = -30085.6
$ws.Range("H122").Value = 5081.364
$ws.Range("I122").Value = 4839.5
$ws.Range("K122").Value = 14518.5
$ws.Range("M122").Value = -12068.5
$ws.Range("H132").Value = 2115.111
$ws.Range("I132").Value = 1276.7142
$ws.Range("J132").Value = 5049.5
$ws.Range("K132").Value = 3830.1426
$ws.Range("L132").Value = 15148.5
$ws.Range("M132").Value = -1300.1426
$ws.Range("N132").Value = -20208.5

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H56").Value = 26096.5
$ws.Range("I56").Value = 26096.5
$ws.Range("K56").Value = 26096.5
$ws.Range("M56").Value = -25566.5
$ws.Range("H75").Value = 1350
$ws.Range("I75").Value = 1350
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 4050
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -3052
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1350
$ws.Range("I78").Value = 1350
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 12150
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -7158
$ws.Range("N78").ClearContents()

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H7").Value = 343467.66
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 100
$ws.Range("M7").Value = 12
$ws.Range("H8").Value = 343467.66
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 100
$ws.Range("M8").Value = 39
$ws.Range("H102").Value = 3018.1904
$ws.Range("I102").Value = 2463.7058
$ws.Range("J102").Value = 5374.75
$ws.Range("K102").Value = 2463.7058
$ws.Range("L102").Value = 5374.75
$ws.Range("M102").Value = -841.7058000000002
$ws.Range("N102").Value = -8618.75
$ws.Range("H113").Value = 1002277.75
$ws.Range("I113").Value = 1336037
$ws.Range("K113").Value = 1336037
$ws.Range("M113").Value = -1333867
$ws.Range("H122").Value = 3383.25
$ws.Range("I122").Value = 1788
$ws.Range("K122").Value = 5364
$ws.Range("M122").Value = -2914
$ws.Range("H126").Value = 4788.75
$ws.Range("I126").Value = 4788.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14366.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11896.25
$ws.Range("N126").ClearContents()

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5280
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H40").Value = 4384.5557
$ws.Range("I40").Value = 4461.2666
$ws.Range("J40").Value = 4001
$ws.Range("K40").Value = 4461.2666
$ws.Range("L40").Value = 4001
$ws.Range("M40").Value = -4325.2666
$ws.Range("N40").Value = -4273
$ws.Range("H55").Value = 1289.8
$ws.Range("J55").Value = 324.5
$ws.Range("L55").Value = 324.5
$ws.Range("N55").Value = -670.5
$ws.Range("H61").Value = 35627.97
$ws.Range("I61").Value = 38873.93
$ws.Range("K61").Value = 38873.93
$ws.Range("M61").Value = -38671.93
$ws.Range("H100").Value = 4549025
$ws.Range("I100").Value = 8336546
$ws.Range("K100").Value = 8336546
$ws.Range("M100").Value = -8336005
$ws.Range("H113").Value = 35627.97
$ws.Range("I113").Value = 38873.93
$ws.Range("K113").Value = 38873.93
$ws.Range("M113").Value = -36703.93

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H13").Value = 3500
$ws.Range("I13").Value = 2000
$ws.Range("K13").Value = 2000
$ws.Range("M13").Value = -1860
$ws.Range("H41").Value = 13369.571
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 13369.571
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 13369.571
$ws.Range("N41").Value = -14149.571
$ws.Range("M41").ClearContents()
$ws.Range("H51").Value = 17333.334
$ws.Range("I51").Value = 16000
$ws.Range("K51").Value = 16000
$ws.Range("M51").Value = -15490
$ws.Range("H113").Value = 1128.5714
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 4350
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2180
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 8609.6
$ws.Range("I126").Value = 5156
$ws.Range("K126").Value = 15468
$ws.Range("M126").Value = -12998
$ws.Range("H136").Value = 4416.6
$ws.Range("J136").Value = 8721.888999999999
$ws.Range("L136").Value = 26165.667
$ws.Range("N136").Value = -31265.667
